$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells with the same style (bold/border/centered) as the
# existing header row (A1 uses style index 1).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record columns for every data row (2-51).
$ws.Range("AD2:AD51").Value = 88
$ws.Range("AE2:AE51").Value = 74
$ws.Range("AF2:AF51").Value = 0
